$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.810.27"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.476.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.29"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.396"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.477.12"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.28"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +11.75%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.123.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "93.658.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000249"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.473.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.44"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.83"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +9.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.526"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +19.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "503.57"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000186"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.642.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.24%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.990"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.561"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.70"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "571.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.63"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.28%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.14%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +10.63%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.75"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.59"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.13"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.99%  "
